$d = $word.ActiveDocument

# =========================================================================
# Part 1: "Dear Nikola Tesla," -> "Dear " / "AAAAAAAAAAAA" / bookmark / ","
# =========================================================================

$greeting = $d.Paragraphs(3).Range
$greeting.Find.Execute("Nikola Tesla", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$nameStart = $greeting.Start
$nameEnd = $greeting.End

# Temporary wall so the new name doesn't re-merge with "Dear ".
$d.Bookmarks.Add("LeftWall1", $d.Range($nameStart, $nameStart))

# Two-step replace (through an intermediate value) forces a real edit so the
# engine doesn't silently keep the old run split/merge state.
$d.Range($nameStart, $nameEnd).Text = "ZZZZZZZZZZZZ"
$lw1 = $d.Bookmarks("LeftWall1")
$d.Range($lw1.End, $lw1.End + 12).Text = "AAAAAAAAAAAA"

$lw1 = $d.Bookmarks("LeftWall1")
$goBackPos = $lw1.End + 12

# Move _GoBack to sit right after the new name (this also walls "AAAAAAAAAAAA"
# off from the trailing comma so they don't re-merge).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))

$d.Bookmarks("LeftWall1").Delete()

# =========================================================================
# Part 2: "I have" + bookmark + " " -> "I have " (single run, no bookmark)
# =========================================================================

$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("I believe ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$leftWallPos = $p4.End

$p4b = $d.Paragraphs(4).Range
$p4b.Find.Execute("a few hundred", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rightWallPos = $p4b.Start

$d.Bookmarks.Add("LeftWall2", $d.Range($leftWallPos, $leftWallPos))
$d.Bookmarks.Add("RightWall2", $d.Range($rightWallPos, $rightWallPos))

# The old _GoBack sits between "I have" and " " -- remove it so they can merge.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Range($leftWallPos, $rightWallPos).Text = "XXXXXXX"
$lw2 = $d.Bookmarks("LeftWall2")
$rw2 = $d.Bookmarks("RightWall2")
$d.Range($lw2.End, $rw2.Start).Text = "I have "

$d.Bookmarks("LeftWall2").Delete()
$d.Bookmarks("RightWall2").Delete()

Write-Output $d.Content.Text
